$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# "Expected Behaviour" cell (H2): append two new Iconposition validation
# lines to the end of the validate4 block, just before its closing "};".
$newText = @"
validate1
{
validate_PageTitle=Manual Compliance Ruby Specs
};
validate2
{
validate_PageTitle=Battery Ruby Test
};
validate3
{
validate_Text_Exists=VT200-0251
};
validate4
{
validate_Screenshot=VT200_0251
validate_Iconposition=batteryview_xpath,left,20
validate_Iconposition=batteryview_xpath,top,40
};
"@

$ws.Range("H2").Value = $newText.Trim()

# The extra two lines make the wrapped, auto-sized cell taller - match
# Excel's recalculated row height for row 2.
$ws.Rows.Item(2).RowHeight = 203.25

# Move the active selection from J2 to G2.
[void]$ws.Range("G2").Select()
